# Updated cryptos list on Fri Feb  2 22:54:08 UTC 2024 with GitHub Actions
# Refreshes the Price (D) / Volume(1h) (E) columns for each coin row, and
# swaps the EnergySwap / NEARProtocol rows (46/47) with their refreshed data.
#
# Note: some refreshed Price values (column D) are plain decimal numbers
# (e.g. "68.00", "0.510"). Excel would normally auto-convert such strings to
# numeric cells (losing the original text formatting / trailing zeros), so
# for those cells we force the cell to Text ("@") before writing the value,
# then restore the cell style to "Normal" so no stray formatting is left
# behind (matching the source workbook, where these cells carry no style).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.127.23"
$ws.Range("E2").Value = "  +0.43%  "
$ws.Range("D3").Value = "2.305.54"
$ws.Range("E3").Value = "  +0.51%  "
$ws.Range("E4").Value = "  +0.03%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "301.13"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.49%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "100.09"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +3.13%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.510"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +1.11%  "
$ws.Range("E9").Value = "  +2.02%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "36.57"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +8.64%  "
$ws.Range("E12").Value = "  +0.58%  "
$ws.Range("E13").Value = "  +4.24%  "
$ws.Range("E14").Value = "  +2.60%  "
$ws.Range("D15").Value = "2.664.80"
$ws.Range("E15").Value = "  +0.65%  "
$ws.Range("D16").Value = "2.309.49"
$ws.Range("E16").Value = "  +0.85%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.799"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -0.69%  "
$ws.Range("D18").Value = "43.014.58"
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "12.82"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +9.98%  "
$ws.Range("E20").Value = "  +0.57%  "
$ws.Range("E21").Value = "  +1.64%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "68.00"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +0.77%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "235.89"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -0.23%  "
$ws.Range("E24").Value = "  +8.54%  "
$ws.Range("E26").Value = "  -0.07%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "25.23"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +3.93%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "170.39"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +2.41%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "34.55"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +2.47%  "
$ws.Range("E30").Value = "  -0.62%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "9.15"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +0.62%  "
$ws.Range("E32").Value = "  +0.00%  "
$ws.Range("E33").Value = "  +2.32%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "17.75"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +5.95%  "
$ws.Range("E35").Value = "  +0.72%  "
$ws.Range("E36").Value = "  -0.37%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.0692"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -0.63%  "
$ws.Range("E38").Value = "  +1.62%  "
$ws.Range("E39").Value = "  +1.70%  "
$ws.Range("E40").Value = "  -0.10%  "
$ws.Range("E41").Value = "  +0.88%  "
$ws.Range("E42").Value = "  +3.12%  "
$ws.Range("D43").Value = "1.987.46"
$ws.Range("E43").Value = "  +0.11%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "2.27"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -4.15%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "10.21"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +3.36%  "
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "2.90"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +2.58%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "17.65"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +0.93%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "55.58"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +4.80%  "
$ws.Range("E49").Value = "  +5.02%  "
$ws.Range("D50").Value = "2.530.23"
$ws.Range("E50").Value = "  -0.17%  "
$ws.Range("E51").Value = "  +1.31%  "
